# Auto-committed on 2022/09/16 週五 17:39:19.43
#
# This edit:
#  1. Adds a new "findByTranNoFirst" lookup row (row 13) to the DBS sheet,
#     reusing the existing "TranNo = ,AND MrKey = " ORDER text and adding a
#     new "TxDate asc" order clause.
#  2. Refreshes the DBD row heights (re-layout after the Excel build bump).
#  3. Leaves the selection on DBD at B10 and switches the active/visible
#     sheet to DBS with the selection parked on the freshly-added A13.

$wb = $excel.ActiveWorkbook
$wsDBD = $wb.Worksheets.Item("DBD")
$wsDBS = $wb.Worksheets.Item("DBS")

# --- DBS: append the new row (A13:C13) ---------------------------------
$wsDBS.Range("A13").Value = "findByTranNoFirst"
$wsDBS.Range("B13").Value = "TranNo = ,AND MrKey = "
$wsDBS.Range("C13").Value = "TxDate asc"

# --- DBD: row-height refresh (matches the re-layout after the Excel
#     build bump) ----------------------------------------------------------
foreach ($r in 1,2,3,4,5,6,9) {
    $wsDBD.Rows.Item($r).AutoFit()
}
$wsDBD.Rows.Item(13).RowHeight = 16.2
$wsDBD.Rows.Item(23).RowHeight = 31.2

# --- Selections / active sheet -----------------------------------------
# DBD keeps focus on B10 (no longer the visible tab) ...
[void]$wsDBD.Activate()
[void]$wsDBD.Range("B10").Select()

# ... while DBS becomes the active/visible tab, selection on the new row.
[void]$wsDBS.Activate()
[void]$wsDBS.Range("A13").Select()
